$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.014.57"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").Value = "2.309.79"
$ws.Range("E3").Value = "  +1.93%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'303.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("E6").Value = "  +5.51%  "
$ws.Range("E7").Value = "  +1.68%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.516"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.90%  "
$ws.Range("D10").Value = "'35.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.79%  "
$ws.Range("E12").Value = "  +4.16%  "
$ws.Range("D13").Value = "'17.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +15.48%  "
$ws.Range("E14").Value = "  +3.61%  "
$ws.Range("D15").Value = "2.686.44"
$ws.Range("E15").Value = "  +2.63%  "
$ws.Range("D16").Value = "2.291.54"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("E17").Value = "  +3.86%  "
$ws.Range("D18").Value = "42.943.54"
$ws.Range("E18").Value = "  +2.06%  "
$ws.Range("E19").Value = "  +8.56%  "
$ws.Range("D20").Value = "'6.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.46%  "
$ws.Range("E21").Value = "  +1.63%  "
$ws.Range("D22").Value = "'67.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.16%  "
$ws.Range("D23").Value = "'237.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.33%  "
$ws.Range("E24").Value = "  +13.61%  "
$ws.Range("E25").Value = "  +1.10%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").Value = "'24.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.59%  "
$ws.Range("E28").Value = "  +5.52%  "
$ws.Range("D29").Value = "'167.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").Value = "'34.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.94%  "
$ws.Range("D31").Value = "'9.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").Value = "'5.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = "  +3.17%  "
$ws.Range("E35").Value = "  +3.60%  "
$ws.Range("D36").Value = "'17.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.20%  "
$ws.Range("E37").Value = "  +0.99%  "
$ws.Range("E38").Value = "  +4.01%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "'2.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'1.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.14%  "
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("D42").Value = "2.004.06"
$ws.Range("E42").Value = "  +2.32%  "
$ws.Range("E43").Value = "  -5.30%  "
$ws.Range("D44").Value = "'0.0288"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.04%  "
$ws.Range("D46").Value = "'17.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("D47").Value = "'2.85"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.55%  "
$ws.Range("D48").Value = "'55.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.54%  "
$ws.Range("D49").Value = "2.529.72"
$ws.Range("E49").Value = "  +1.55%  "
$ws.Range("D50").Value = "'1.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.09%  "
